# Daily update at 8 AM UTC
# Adds the next day's row (row 66) to the "Wins Over Time" sheet and
# moves the "most recent row" date-only formatting from the old last
# row (65) to the new last row (66).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 65 is no longer the most recent entry, so it goes back to the
# regular date/time number format used by every other data row.
$ws.Range("A65").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data in row 66.
$ws.Range("A66").Value = 45653
$ws.Range("B66").Value = 156
$ws.Range("C66").Value = 146
$ws.Range("D66").Value = 153

# Row 66 is now the most recent entry, so its date cell gets the
# date-only "highlight" number format.
$ws.Range("A66").NumberFormat = "YYYY-MM-DD"
